# Updated cryptos list values (prices and 1h volume change %)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.495.51"
$ws.Range("E2").Value = "  +1.84%  "
$ws.Range("D3").Value = "2.284.18"
$ws.Range("E3").Value = "  +0.91%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.66"
$ws.Range("E5").Value = "  +1.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.25"
$ws.Range("E6").Value = "  +6.92%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.531"
$ws.Range("E7").Value = "  +0.46%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +3.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.25"
$ws.Range("E10").Value = "  +12.32%  "
$ws.Range("E11").Value = "  +0.95%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.112"
$ws.Range("E12").Value = "  -2.06%  "
$ws.Range("E13").Value = "  +1.69%  "
$ws.Range("D14").Value = "2.638.46"
$ws.Range("E14").Value = "  +0.93%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.54"
$ws.Range("E15").Value = "  +2.39%  "
$ws.Range("D16").Value = "2.285.06"
$ws.Range("E16").Value = "  +1.50%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.804"
$ws.Range("E17").Value = "  +5.31%  "
$ws.Range("D18").Value = "42.363.59"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.67"
$ws.Range("E19").Value = "  +0.82%  "
$ws.Range("E20").Value = "  +1.56%  "
$ws.Range("E21").Value = "  +1.91%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.98"
$ws.Range("E22").Value = "  +1.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "242.65"
$ws.Range("E23").Value = "  +1.26%  "
$ws.Range("E24").Value = "  +0.87%  "
$ws.Range("E25").Value = "  +1.96%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "38.16"
$ws.Range("E28").Value = "  +10.82%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.60"
$ws.Range("E29").Value = "  +0.72%  "
$ws.Range("E30").Value = "  +2.89%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "161.11"
$ws.Range("E31").Value = "  +0.76%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.28"
$ws.Range("E32").Value = "  +0.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  +0.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.15"
$ws.Range("E34").Value = "  +5.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0750"
$ws.Range("E35").Value = "  +0.89%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.29"
$ws.Range("E36").Value = "  +2.74%  "
$ws.Range("E37").Value = "  +2.78%  "
$ws.Range("E38").Value = "  +0.59%  "
$ws.Range("E39").Value = "  +3.90%  "
$ws.Range("E40").Value = "  -0.67%  "
$ws.Range("E41").Value = "  +6.34%  "
$ws.Range("E42").Value = "  +13.82%  "
$ws.Range("D43").Value = "2.004.61"
$ws.Range("E43").Value = "  -1.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.35"
$ws.Range("E44").Value = "  +0.64%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0287"
$ws.Range("E45").Value = "  +3.32%  "
$ws.Range("E46").Value = "  +4.88%  "
$ws.Range("E47").Value = "  -2.62%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "53.63"
$ws.Range("E48").Value = "  +3.35%  "
$ws.Range("E49").Value = "  +1.81%  "

# Row 50/51: Aave and BitcoinSV swap positions with updated values
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "93.53"
$ws.Range("E50").Value = "  +2.67%  "
$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.68"
$ws.Range("E51").Value = "  -0.11%  "
